$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Abfrage: {hw_10" + "_file}" were stored as two separate runs (an
#    artifact of earlier editing). Re-assert the same text via Find/Replace
#    so Word collapses it back into a single run -- matches runs hw_1..hw_9
#    which already live in a single run each.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Abfrage: {hw_10_file}", $true, $false, $false, $false, $false,
    $true, 1, $false, "Abfrage: {hw_10_file}", 2)

# ---------------------------------------------------------------------------
# 2) "{@hw_10" + bookmark(_GoBack) + "}" likewise need their two text runs
#    merged into one run reading "{@hw_10}", while keeping the _GoBack
#    bookmark (it simply ends up anchored after the merged run instead of
#    in the middle of the old split). We do this surgically so the
#    bookmark is never dropped:
#      a) delete the lone "}" run's character
#      b) insert "}" back in right before the (now-trailing) bookmark
# ---------------------------------------------------------------------------
$hw10Para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "{@hw_10}") {
        $hw10Para = $p
        break
    }
}

if ($hw10Para -ne $null) {
    $closeBracePos = $hw10Para.Range.Start + 7   # offset of the "}" character
    $braceRange = $d.Range($closeBracePos, $closeBracePos + 1)
    if ($braceRange.Text -eq "}") {
        $braceRange.Delete()
        $insertPos = $d.Range($closeBracePos, $closeBracePos)
        $insertPos.InsertBefore("}")
    }
}

# ---------------------------------------------------------------------------
# 3) Drop the five trailing empty paragraphs after the {@hw_10} paragraph
#    (right before the section break) -- delete the paragraph mark at the
#    end of that paragraph, one at a time, so it keeps absorbing the next
#    (empty) paragraph until none of the blank ones remain.
# ---------------------------------------------------------------------------
for ($n = 0; $n -lt 5; $n++) {
    $hw10Para = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13) -eq "{@hw_10}") {
            $hw10Para = $p
        }
    }
    if ($hw10Para -eq $null) { break }

    $nextPara = $hw10Para.Next()
    if ($nextPara -eq $null) { break }
    if ($nextPara.Range.Text.TrimEnd([char]13) -ne "") { break }

    $markPos = $hw10Para.Range.End - 1
    $markRange = $d.Range($markPos, $markPos + 1)
    $markRange.Delete()
}
